# Weekly update: insert two new daily-price records for "Cilantro" at
# Terminal Hortofrutícola Agro Chillán, ahead of the existing rows 47-60,
# pushing the rest of the table down by two rows (table grows from 60 to
# 62 data/rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 47 (shifts rows 47:60 down to 49:62,
# inheriting formatting/number-format from the row above, same as Excel's
# default "Insert" behavior).
$ws.Range("A47:A48").EntireRow.Insert()

# New row 47
$ws.Cells.Item(47, 1).Value2 = 7
$ws.Cells.Item(47, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(47, 3).Value2 = "Ñuble"
$ws.Cells.Item(47, 4).Value2 = 44767
$ws.Cells.Item(47, 5).Value2 = 16
$ws.Cells.Item(47, 6).Value2 = 100112040
$ws.Cells.Item(47, 7).Value2 = "Cilantro"
$ws.Cells.Item(47, 8).Value2 = "Sin especificar"
$ws.Cells.Item(47, 9).Value2 = "Primera"
$ws.Cells.Item(47, 10).Value2 = 200
$ws.Cells.Item(47, 11).Value2 = 700
$ws.Cells.Item(47, 12).Value2 = 800
$ws.Cells.Item(47, 13).Value2 = 750
$ws.Cells.Item(47, 14).Value2 = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(47, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(47, 16).Value2 = 750
$ws.Cells.Item(47, 17).Value2 = 1
$ws.Cells.Item(47, 18).Value2 = "Hortaliza"

# New row 48
$ws.Cells.Item(48, 1).Value2 = 7
$ws.Cells.Item(48, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(48, 3).Value2 = "Ñuble"
$ws.Cells.Item(48, 4).Value2 = 44767
$ws.Cells.Item(48, 5).Value2 = 16
$ws.Cells.Item(48, 6).Value2 = 100112040
$ws.Cells.Item(48, 7).Value2 = "Cilantro"
$ws.Cells.Item(48, 8).Value2 = "Sin especificar"
$ws.Cells.Item(48, 9).Value2 = "Segunda"
$ws.Cells.Item(48, 10).Value2 = 150
$ws.Cells.Item(48, 11).Value2 = 600
$ws.Cells.Item(48, 12).Value2 = 600
$ws.Cells.Item(48, 13).Value2 = 600
$ws.Cells.Item(48, 14).Value2 = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(48, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(48, 16).Value2 = 600
$ws.Cells.Item(48, 17).Value2 = 1
$ws.Cells.Item(48, 18).Value2 = "Hortaliza"
